$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Лабораторна робота №2" -> "Лабораторна робота №3"
#    Scope the Find to the title paragraph only, and match just the
#    trailing digit, so only that single run's text changes (the
#    sibling runs "Лабораторна" / " робота №" stay untouched, exactly
#    like the source diff).
# ---------------------------------------------------------------------
$titleRange = $d.Paragraphs(8).Range
$titleRange.Find.Execute("2", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "3", 2)

# ---------------------------------------------------------------------
# 2. Theme paragraph: "«Використання програмних пакетів для розробки
#    CSS-оформлення.»" -> "«Використання програмних пакетів для
#    інтернет-графіки»", ending up split across 5 runs with the
#    "_GoBack" bookmark relocated in between "Викорис" and
#    "тання програмн".
# ---------------------------------------------------------------------
$themeRange = $d.Paragraphs(10).Range
$themeRange.Find.Execute("«Використання програмних пакетів для розробки CSS-оформлення.»", `
                          $false, $false, $false, $false, $false, `
                          $true, 1, $false, `
                          "«Використання програмних пакетів для інтернет-графіки»", 2)

# Locate the freshly-written sentence so the split points are computed
# from its actual (current) position rather than a hard-coded offset.
$find = $d.Paragraphs(10).Range.Duplicate
$find.Find.Execute("«Використання програмних пакетів для інтернет-графіки»")
$base = $find.Start

$p1 = "«"
$p2 = "Викорис"
$p3 = "тання програмн"
$p4 = "их пакетів для інтернет-графіки"
$p5 = "»"

$s1 = $base
$s2 = $s1 + $p1.Length
$s3 = $s2 + $p2.Length
$s4 = $s3 + $p3.Length
$s5 = $s4 + $p4.Length
$s6 = $s5 + $p5.Length

# Force run boundaries at each interior split point by toggling a
# character format on/off across the "Викорис" and "тання програмн ...
# графіки" chunks -- Word must break the run apart to store the
# (temporary) differing formatting, and once the value is restored the
# pieces remain separate runs, matching the target run layout.
$chunkB = $d.Range($s2, $s3)
$chunkB.Font.Bold = 1
$chunkB2 = $d.Range($s2, $s3)
$chunkB2.Font.Bold = 0

$chunkD = $d.Range($s4, $s5)
$chunkD.Font.Bold = 1
$chunkD2 = $d.Range($s4, $s5)
$chunkD2.Font.Bold = 0

# Move the "_GoBack" bookmark (auto-tracks "last edit") to sit right
# between the "Викорис" and "тання програмн" runs -- adding it here
# removes it from wherever it previously lived (the empty paragraph
# further down).
$goBackRange = $d.Range($s3, $s3)
$d.Bookmarks.Add("_GoBack", $goBackRange)
